$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.778.11"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "3.461.16"
$ws.Range("E3").Value = "  +2.47%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "576.56"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.59%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "148.64"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +8.38%  "
$ws.Range("D7").Value = "3.461.30"
$ws.Range("E7").Value = "  +2.39%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  +0.23%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "7.71"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +3.47%  "
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").Value = "4.046.59"
$ws.Range("E13").Value = "  +2.48%  "
$ws.Range("E14").Value = "  -1.30%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "27.24"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +4.61%  "
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").Value = "3.447.29"
$ws.Range("E17").Value = "  +2.38%  "
$ws.Range("D18").Value = "61.832.35"
$ws.Range("E18").Value = "  +0.39%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.16"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +4.64%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "14.14"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.36%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "9.58"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.36%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "383.35"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.41%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.563"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("D24").Value = "3.573.77"
$ws.Range("E24").Value = "  +1.75%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.14%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "72.36"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.32%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.0000126"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.32%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.178"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +9.28%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "7.83"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +3.82%  "
$ws.Range("E30").Value = "  -9.29%  "
$ws.Range("E31").Value = "  +0.04%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "8.23"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("E34").Value = "  -0.03%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "24.00"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +1.67%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "5.30"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.34%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "7.02"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.53%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.58"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +2.94%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "166.92"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.31%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.0791"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +2.61%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "26.28"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +7.69%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.796"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.87%  "
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("E44").Value = "  +0.34%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "42.28"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("E46").Value = "  +1.90%  "
$ws.Range("E47").Value = "  -0.35%  "
$ws.Range("D48").Value = "2.649.77"
$ws.Range("E48").Value = "  +12.62%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "23.92"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +4.85%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "6.91"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.93%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.20"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +8.03%  "
